# Add a new row (row 30) to each of the 4 sheets, duplicating row 29's
# data but with an updated timestamp in column A (next day's reading).
# This mirrors a CSV-ingest appending one more sample per sheet
# (commit: "Add csv module error handling").

$wb = $excel.ActiveWorkbook

$newTimestamp = 45816.49134259259

$rowsData = @{
    "FE_LFT_#1" = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x74"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 372
        I = 15
    }
    "FE_LFT_#2" = @{
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x88"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 392
        I = 14
    }
    "FE_PLT_#1" = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x6D"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 109
        I = 3
    }
    "FE_PLT_#2" = @{
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x6D"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 109
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $rowsData.ContainsKey($name)) { continue }
    $data = $rowsData[$name]

    $newRow = 30

    $ws.Cells.Item($newRow, 1).Value = $newTimestamp
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}

Write-Output "Row 30 added to all sheets"
